# "New future enhancement added"
# Adds three new "To Do" rows to Sheet1 (rows 6-8), bolds the header row,
# widens the Module column to fit the new, longer entries, and sets the
# page orientation to portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bold the header row (A1:F1) -- introduces the new bold font / cellXf.
$ws.Range("A1:F1").Font.Bold = $true

# 2) New rows of data.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Multiple Payment for one invoice"
$ws.Range("C6").Value = "Invoice"
$ws.Range("D6").Value = "Open"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Revert check out in case of operator mistake"
$ws.Range("C7").Value = "Checkin"
$ws.Range("D7").Value = "Open"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Recycle bin"
$ws.Range("C8").Value = "Navigator and all component"
$ws.Range("D8").Value = "Open"

# 3) Widen column C (Module) so the new, longer entries fit.
$ws.Columns("C").ColumnWidth = 26.1666666666667

# 4) Page orientation -> portrait.
$ws.PageSetup.Orientation = 1

# 5) Leave the active cell on the last entry, like the original author did.
$ws.Range("B8").Select() | Out-Null
